# Commit: "commiting POM framework for selenium"
# Adds a block of repeated login rows (rows 3-11, mirroring row 2) to the
# existing LoginTestData sheet, each with a hyperlinked password cell, and
# appends a new blank "Sheet1" worksheet after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill rows 3..11 with the same Username/Password pair as row 2, each
# Password cell carrying its own hyperlink (mirrors B2's mailto link).
for ($r = 3; $r -le 11; $r++) {
    $ws.Range("A$r").Value = 9594748758
    $ws.Range("B$r").Value = "Mokalpur@2021"
    $ws.Hyperlinks.Add($ws.Range("B$r"), "mailto:Mokalpur@2021") | Out-Null
    $ws.Range("B$r").Style = "Hyperlink"
}

# Move the selection to B2, as in the saved workbook.
$ws.Range("B2").Select() | Out-Null

# Append a new empty worksheet ("Sheet1") after LoginTestData, then restore
# LoginTestData as the active/selected sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet) | Out-Null
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
